# Add a new row of calibrator results ("Add more result other pc")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "traffic_factor"
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 48
$ws.Range("D12").Value = 7
$ws.Range("E12").Value = 0.43
$ws.Range("F12").Value = 0.08

$ws.Range("E13").Select()
